$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-18 Sunday" "2026-01-19 Monday"

Replace-Text "26×38=988" "61×89=5429"
Replace-Text "48×99=4752" "36×53=1908"
Replace-Text "32×98=3136" "32×16=512"
Replace-Text "54×97=5238" "80×72=5760"
Replace-Text "90×89=8010" "79×23=1817"
Replace-Text "49×39=1911" "76×79=6004"
Replace-Text "64×71=4544" "62×71=4402"
Replace-Text "56×71=3976" "31×31=961"
Replace-Text "48×27=1296" "71×92=6532"
Replace-Text "17×78=1326" "81×62=5022"
Replace-Text "64×95=6080" "97×49=4753"
Replace-Text "80×14=1120" "82×16=1312"
Replace-Text "90×84=7560" "15×94=1410"
Replace-Text "38×31=1178" "25×21=525"
Replace-Text "50×61=3050" "34×54=1836"
Replace-Text "75×90=6750" "33×62=2046"
Replace-Text "67×69=4623" "19×47=893"
Replace-Text "12×27=324" "13×82=1066"
Replace-Text "99×50=4950" "78×54=4212"
Replace-Text "86×71=6106" "94×75=7050"
Replace-Text "70×16=1120" "78×51=3978"
Replace-Text "26×47=1222" "93×90=8370"
Replace-Text "96×84=8064" "66×61=4026"
Replace-Text "22×21=462" "13×76=988"
Replace-Text "61×63=3843" "13×25=325"
